# Update gh-pages output data (苏州-漫展信息) — refresh scraped counters.
#
# Sheet "展览" (Exhibitions) and its mirror "全部类型" (All types) share most
# rows; "演出" (Performances) feeds the remaining rows of "全部类型".
# Each sheet is patched independently since the workbook stores the merged
# "全部类型" view as its own physical rows, not a formula/reference.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 ("F" = want-to-go count, "G" = lowest price, "E" = time range) ----
$ws1.Range("F3").Value = 1063
$ws1.Range("F4").Value = 536
$ws1.Range("E5").Value = "2024.09.15 10:00-09.17 17:00"
$ws1.Range("F5").Value = 13977
$ws1.Range("G6").Value = "不可售"
$ws1.Range("F7").Value = 569
$ws1.Range("F8").Value = 229
$ws1.Range("F9").Value = 1802
$ws1.Range("F11").Value = 147
$ws1.Range("F14").Value = 545
$ws1.Range("F17").Value = 3
$ws1.Range("F18").Value = 14068
$ws1.Range("F19").Value = 373
$ws1.Range("F21").Value = 15008
$ws1.Range("F22").Value = 15
$ws1.Range("F23").Value = 8333
$ws1.Range("F24").Value = 283
$ws1.Range("F26").Value = 34
$ws1.Range("F27").Value = 158
$ws1.Range("F29").Value = 166
$ws1.Range("F32").Value = 3
$ws1.Range("F33").Value = 32
$ws1.Range("F34").Value = 1043
$ws1.Range("F35").Value = 26
$ws1.Range("F40").Value = 12
$ws1.Range("F41").Value = 219
$ws1.Range("F42").Value = 225
$ws1.Range("F43").Value = 395
$ws1.Range("F45").Value = 5129

# ---- 演出 ----
$ws2.Range("F2").Value = 52

# ---- 全部类型 (merged view: rows 2-37 mirror 展览 rows 2-37, rows 38-39
#      mirror 演出 rows 2-3, rows 40-48 mirror 展览 rows 38-46) ----
$ws4.Range("F3").Value = 1063
$ws4.Range("F4").Value = 536
$ws4.Range("E5").Value = "2024.09.15 10:00-09.17 17:00"
$ws4.Range("F5").Value = 13977
$ws4.Range("G6").Value = "不可售"
$ws4.Range("F7").Value = 569
$ws4.Range("F8").Value = 229
$ws4.Range("F9").Value = 1802
$ws4.Range("F11").Value = 147
$ws4.Range("F14").Value = 545
$ws4.Range("F17").Value = 3
$ws4.Range("F18").Value = 14068
$ws4.Range("F19").Value = 373
$ws4.Range("F21").Value = 15008
$ws4.Range("F22").Value = 15
$ws4.Range("F23").Value = 8333
$ws4.Range("F24").Value = 284
$ws4.Range("F26").Value = 34
$ws4.Range("F27").Value = 158
$ws4.Range("F29").Value = 166
$ws4.Range("F32").Value = 3
$ws4.Range("F33").Value = 32
$ws4.Range("F34").Value = 1043
$ws4.Range("F35").Value = 26
$ws4.Range("F38").Value = 52
$ws4.Range("F42").Value = 12
$ws4.Range("F43").Value = 219
$ws4.Range("F44").Value = 225
$ws4.Range("F45").Value = 395
$ws4.Range("F47").Value = 5129

Write-Output "applied"
